# Updates cryptos list price (D) and volume-change (E) columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.305.25"
$ws.Range("E2").Value = "  +0.29%  "
# Row 3
$ws.Range("D3").Value = "1.876.18"
$ws.Range("E3").Value = "  +0.88%  "
# Row 4
$ws.Range("E4").Value = "  +0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7121"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.65%  "
# Row 7
$ws.Range("E7").Value = "  -0.01%  "
# Row 8
$ws.Range("E8").Value = "  +0.86%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07760"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.11%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.86%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08491"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.82%  "
# Row 12
$ws.Range("D12").Value = "1.900.06"
$ws.Range("E12").Value = "  +2.35%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.210"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.71%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7100"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.18%  "
# Row 15
$ws.Range("E15").Value = "  +1.17%  "
# Row 16
$ws.Range("D16").Value = "29.306.52"
$ws.Range("E16").Value = "  +0.42%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008245"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.74%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.010"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.35%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.90%  "
# Row 20
$ws.Range("D20").Value = "2.133.33"
$ws.Range("E20").Value = "  +1.84%  "
# Row 21
$ws.Range("E21").Value = "  +0.66%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.812"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.24%  "
# Row 24
$ws.Range("E24").Value = "  -0.02%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1618"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.12%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.019"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.02%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.07%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.513"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.01%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.400"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.305"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.32%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.277"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.97%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05238"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.84%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.931"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.01%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.176"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.26%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7412"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.87%  "
# Row 37
$ws.Range("E37").Value = "  +0.37%  "
# Row 38
$ws.Range("E38").Value = "  +0.46%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.724"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.44%  "
# Row 40
$ws.Range("D40").Value = "1.176.86"
$ws.Range("E40").Value = "  +2.15%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.378"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.73%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8885"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.82%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.82%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.62%  "
# Row 45
$ws.Range("E45").Value = "  +0.02%  "
# Row 46
$ws.Range("D46").Value = "2.029.63"
$ws.Range("E46").Value = "  +1.46%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.810"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.35%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5209"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.43%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000121"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.47%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.384"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.54%  "
# Row 51
$ws.Range("E51").Value = "  +0.86%  "
